# Update column F ("dSF") values for specific rows on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    8  = -4
    9  = -4
    10 = 4
    16 = 7
    20 = 1
    22 = 1
    23 = -5
    24 = -7
    30 = 2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
